# edit.ps1
# Applies the CasosColombia.xlsx update:
#  - Patches a batch of existing cells that flip between numeric values
#    and the literal "NaN" placeholder (shared string) across rows
#    12-198 (data corrections in the cumulative/deceased columns).
#  - Appends a brand-new data row (209) for 2020-09-29 (serial 44103).
#  - Moves the active selection to the new row's first cell, matching
#    the source workbook's saved cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CasosColombia")

$patches = @{
    "BU12" = "NaN"
    "W36" = "NaN"
    "W43" = 9
    "BQ88" = "NaN"
    "DG104" = "NaN"
    "CG114" = 25
    "CG115" = 30
    "CG120" = 42
    "CG121" = 42
    "CG126" = 72
    "CG127" = 73
    "CG130" = 100
    "CG136" = 195
    "CG137" = 296
    "CG139" = 415
    "CG140" = 439
    "CG141" = 447
    "CG142" = 451
    "CG143" = 504
    "CG144" = 547
    "CG145" = 626
    "AG146" = 14
    "CG146" = 681
    "CG147" = 747
    "CG148" = 806
    "CG149" = 851
    "CG150" = 1009
    "AH151" = 2
    "CG151" = 1150
    "AH152" = 2
    "CG152" = 1310
    "AH153" = 2
    "CG153" = 1432
    "CG154" = 1581
    "AH155" = 5
    "CG155" = 1609
    "AG156" = 14
    "CG156" = 1770
    "CG157" = 1861
    "CG158" = 2147
    "CG159" = 2215
    "AG160" = 14
    "CG160" = 2338
    "AG161" = 14
    "CG161" = 2462
    "AG162" = 14
    "CG162" = 2538
    "CG163" = 2690
    "CG164" = 2848
    "AG165" = 16
    "CG165" = 3033
    "AG166" = 17
    "CG166" = 3141
    "AG167" = 21
    "CG167" = 3274
    "AG168" = 21
    "CG168" = 3445
    "CG169" = 3520
    "CG170" = 3656
    "CG171" = 3713
    "CG172" = 3806
    "AG173" = 49
    "CG173" = 3997
    "CG174" = 4088
    "CG175" = 4111
    "CG176" = 4214
    "CG177" = 4331
    "CG178" = 4467
    "CG179" = 4575
    "CG180" = 4724
    "CG181" = 4800
    "CG182" = 4920
    "CG183" = 4924
    "CG184" = 5048
    "CG185" = 5115
    "CG186" = 5238
    "CG187" = 5322
    "CG188" = 5428
    "CG189" = 5492
    "CG190" = 5541
    "CG191" = 5580
    "CG192" = 5587
    "CG193" = 5631
    "CG194" = 5641
    "CG195" = 5666
    "CG196" = 5677
    "CG197" = 5701
    "CG198" = 5810
}

# --- Cell-level corrections on existing rows -----------------------------
foreach ($ref in $patches.Keys) {
    $ws.Range($ref).Value = $patches[$ref]
}


# --- New row 209 (2020-09-29) ---------------------------------------------
$row209 = @{
    "B" = 824042
    "C" = 2742
    "D" = 113776
    "E" = 67496
    "F" = 267166
    "G" = 28810
    "H" = 7026
    "I" = 5867
    "J" = 8476
    "K" = 9403
    "L" = 19795
    "M" = 3996
    "N" = 23943
    "O" = 33105
    "P" = 8101
    "Q" = 11389
    "R" = 15320
    "S" = 15139
    "T" = 18321
    "U" = 15691
    "V" = 3771
    "W" = 3541
    "X" = 10906
    "Y" = 30927
    "Z" = 14071
    "AA" = 12227
    "AB" = 62119
    "AC" = 2354
    "AD" = 1424
    "AE" = 765
    "AF" = 476
    "AG" = 869
    "AH" = 527
    "AI" = 767
    "AJ" = 2069
    "AK" = 5869
    "AL" = 38137
    "AN" = 2571
    "AO" = 47918
    "AP" = 1112
    "AQ" = 22973
    "AR" = 1534
    "AS" = 10591
    "AT" = 1679
    "AU" = 1609
    "AV" = 8552
    "AW" = 2035
    "AX" = 964
    "AY" = 2501
    "AZ" = 2692
    "BA" = 65868
    "BB" = 14211
    "BC" = 6789
    "BD" = 9933
    "BE" = 7338
    "BF" = 257
    "BG" = 1477
    "BH" = 2736
    "BI" = 744
    "BJ" = 2169
    "BK" = 9998
    "BL" = 9599
    "BM" = 10751
    "BN" = 14362
    "BO" = 1970
    "BP" = 907
    "BQ" = 14127
    "BR" = 11611
    "BS" = 13521
    "BT" = 3215
    "BU" = 2289
    "BV" = 6120
    "BW" = 4975
    "BX" = 2549
    "BY" = 6125
    "BZ" = 3856
    "CA" = 2259
    "CB" = 1248
    "CC" = 3081
    "CD" = 2249
    "CE" = 2079
    "CF" = 1912
    "CG" = 6605
    "CH" = 2259
    "CI" = 1505
    "CJ" = 1842
    "CK" = 2140
    "CL" = 2275
    "CM" = 2658
    "CN" = 1827
    "CO" = 1227
    "CP" = 1227
    "CQ" = 1091
    "CR" = 3471
    "CS" = 1507
    "CT" = 980
    "CU" = 1145
    "CV" = 1761
    "CW" = 1629
    "CX" = 815
    "CY" = 907
    "CZ" = 1373
    "DA" = 1713
    "DB" = 1677
    "DC" = 1657
    "DD" = 1271
    "DE" = 334
    "DF" = 372
    "DG" = 838
    "DH" = 786
    "DI" = 500
    "DJ" = 544
    "DK" = 389
    "DL" = 674
    "DM" = 760
    "DN" = 527
    "DO" = 507
    "DP" = 374
    "DQ" = 527
    "DR" = 139522
    "DS" = 347773
    "DT" = 20587
    "DU" = 151780
    "DV" = 93839
    "DW" = 46667
    "DX" = 13116
}

$newRow = 209
$ws.Range("A" + $newRow).Value = 44103
foreach ($col in $row209.Keys) {
    $ws.Range($col + $newRow).Value = $row209[$col]
}
$ws.Range("AM" + $newRow).Value = "####"


# --- Restore the saved selection on the new row ---------------------------
$ws.Activate()
$ws.Range("A209").Select()
